$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reverse-Array($arr) {
    if ($arr.Count -le 1) { return $arr }
    return $arr[($arr.Count - 1)..0]
}

function Transform-RecordedBy($s) {
    if ($s -eq $null) { return $s }
    $parts = @()
    foreach ($p in $s.Split(",")) { $parts += $p.Trim() }
    if ($parts.Count -le 1) { return $s }

    if ($parts[0] -eq "system") {
        # Keep a leading lowercase "system" token fixed in place, reverse the rest.
        $rest = $parts[1..($parts.Count - 1)]
        $rest = Reverse-Array $rest
        $newParts = @($parts[0]) + $rest
    } else {
        $newParts = Reverse-Array $parts
    }

    return ($newParts -join ", ")
}

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $null) { continue }
    if ($val -isnot [string]) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $newVal = Transform-RecordedBy $val
    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
